$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.984.29'
$ws.Range("E2").Value = '  +4.85%  '
$ws.Range("D3").Value = '''1.782.29'
$ws.Range("E3").Value = '  +3.53%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = '''243.82'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '''0.4896'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").Value = '''0.2679'
$ws.Range("E8").Value = '  +2.57%  '
$ws.Range("D9").Value = '''0.06257'
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").Value = '''1.781.39'
$ws.Range("E10").Value = '  +3.45%  '
$ws.Range("D11").Value = '''16.39'
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("D12").Value = '''0.07002'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '''0.6269'
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("D14").Value = '''4.647'
$ws.Range("E14").Value = '  +3.25%  '
$ws.Range("D15").Value = '''79.87'
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").Value = '''27.965.05'
$ws.Range("D17").Value = '''1.000'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '''0.9991'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = '''0.000007215'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '''11.95'
$ws.Range("E20").Value = '  +4.84%  '
$ws.Range("D21").Value = '''2.007.50'
$ws.Range("E21").Value = '  +3.14%  '
$ws.Range("D22").Value = '''4.585'
$ws.Range("E22").Value = '  +3.52%  '
$ws.Range("D23").Value = '''8.701'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").Value = '''5.227'
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("D25").Value = '''141.71'
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("D26").Value = '''15.73'
$ws.Range("E26").Value = '  +2.48%  '
$ws.Range("D27").Value = '''1.862'
$ws.Range("E27").Value = '  +6.80%  '
$ws.Range("D28").Value = '''109.32'
$ws.Range("E28").Value = '  +3.12%  '
$ws.Range("D29").Value = '''1.390'
$ws.Range("E29").Value = '  -2.57%  '
$ws.Range("D30").Value = '''4.198'
$ws.Range("E30").Value = '  +6.88%  '
$ws.Range("D31").Value = '''0.08270'
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("D32").Value = '''3.788'
$ws.Range("E32").Value = '  +3.39%  '
$ws.Range("D33").Value = '''0.04799'
$ws.Range("E33").Value = '  +6.53%  '
$ws.Range("D34").Value = '''1.075'
$ws.Range("E34").Value = '  +7.20%  '
$ws.Range("D35").Value = '''2.616'
$ws.Range("E35").Value = '  +0.29%  '
$ws.Range("D36").Value = '''0.6445'
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = '''0.9463'
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("D38").Value = '''2.595'
$ws.Range("E38").Value = '  +7.42%  '
$ws.Range("D39").Value = '''2.053'
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").Value = '''5.928'
$ws.Range("E40").Value = '  +6.31%  '
$ws.Range("D41").Value = '''0.01543'
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("D42").Value = '''0.9995'
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").Value = '''99.97'
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("D44").Value = '''0.3980'
$ws.Range("E44").Value = '  +3.33%  '
$ws.Range("D45").Value = '''7.252'
$ws.Range("E45").Value = '  +4.83%  '
$ws.Range("D46").Value = '''0.1199'
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("D47").Value = '''0.05421'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '''7.981'
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = '''1.294'
$ws.Range("E49").Value = '  +5.12%  '
$ws.Range("D50").Value = '''30.60'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("D51").Value = '''52.83'
$ws.Range("E51").Value = '  +2.25%  '
